$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price updates for rows 2-26 (column D) ---
$priceUpdates = @{
    2  = "247.84"
    3  = "22.48"
    4  = "5.387"
    5  = "0.05673"
    6  = "3.402"
    7  = "6.317"
    8  = "0.8059"
    9  = "0.9231"
    10 = "0.1403"
    11 = "0.07439"
    12 = "0.03085"
    13 = "0.03021"
    14 = "0.09366"
    15 = "3.787"
    16 = "0.001588"
    17 = "0.04737"
    18 = "0.01828"
    19 = "0.0005854"
    20 = "0.006455"
    21 = "0.004974"
    22 = "0.001007"
    23 = "0.0001500"
    24 = "3.690"
    25 = "2.163"
    26 = "0.3254"
}

foreach ($row in $priceUpdates.Keys) {
    $cell = $ws.Range("D$row")
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$row]
}

# --- Rows 41-43: coins rotated (41 <- old42, 42 <- old43, 43 <- old41), with slightly updated prices ---
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1065"
$ws.Range("E41").Value = "40BKEXTokenBKK"

$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.002710"
$ws.Range("E42").Value = "41CEJICEJI"

$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003010"
$ws.Range("E43").Value = "42KickTokenKICKWorstin24h"

# --- Remaining price updates ---
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007525"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005805"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.2091"
